# Update Data by bot, scripted by HH
# Applies the row-2 data refresh (report period rolled back from
# 2020-06-30 to 2019-12-31, DATE_TYPE_CODE 002 -> 001, and the
# associated financial figures / ratios) to Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# DATE_TYPE_CODE is a text code that looks numeric ("001"); use a leading
# apostrophe so Excel keeps it as text instead of coercing it to 1.
$ws.Range("J2").Value = "'001"

# REPORT_DATE is stored as plain text "yyyy-MM-dd 00:00:00", not a real date.
$ws.Range("N2").Value = "2019-12-31 00:00:00"

# Numeric financial figures / ratios for row 2.
$ws.Range("O2").Value = 484611735.47
$ws.Range("P2").Value = 7813887.14
$ws.Range("Q2").Value = 179135841.16
$ws.Range("R2").Value = 51.5039545235
$ws.Range("S2").Value = 205440623.29
$ws.Range("T2").Value = 35.4360445262
$ws.Range("U2").Value = 3051679.69
$ws.Range("V2").Value = 31.3192400061
$ws.Range("W2").Value = 217256873.96
$ws.Range("X2").Value = 57228409.38
$ws.Range("Y2").Value = 23.097634039
$ws.Range("Z2").Value = 10162124.65
$ws.Range("AA2").Value = 55.7251334259
$ws.Range("AB2").Value = 267354861.51
$ws.Range("AC2").Value = 40.0494440561
$ws.Range("AD2").Value = 32.5307579426
$ws.Range("AE2").Value = 24.3176431543
$ws.Range("AF2").Value = 191.0748095619
$ws.Range("AG2").Value = 44.8311210931
